# Add a new "Slots" column (E) to Sheet1, with per-row slot status values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "Slots"
$ws.Range("E2").Value = "Slot E"
$ws.Range("E3").Value = "No Slots"
$ws.Range("E4").Value = "No Slots"
$ws.Range("E5").Value = "No Slots"
$ws.Range("E6").Value = "No Slots"
$ws.Range("E7").Value = "No Slots"
